# Automatic update of files.
# Updates the "Förändrad" (Changed) date column (C) for all data rows
# on the active sheet from 45174 (2023-09-05) to 45175 (2023-09-06).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 66
$col = 3  # Column C = "Förändrad"

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, $col).Value = 45175
}
